$d = $word.ActiveDocument

# --- Text replacements -------------------------------------------------

# "LONG" title line gets " UH" appended.
$d.Content.Find.Execute("LONG", $true, $false, $false, $false, $false, $true, 1, $false, "LONG UH", 2) | Out-Null

# Placeholder "ADDRESS" becomes "HILO".
$d.Content.Find.Execute("ADDRESS", $true, $false, $false, $false, $false, $true, 1, $false, "HILO", 2) | Out-Null

# Placeholder "SECTIONNO" becomes "169".
$d.Content.Find.Execute("SECTIONNO", $true, $false, $false, $false, $false, $true, 1, $false, "169", 2) | Out-Null

# Placeholder "DATE" becomes an actual date.
$d.Content.Find.Execute("DATE", $true, $false, $false, $false, $false, $true, 1, $false, "09/10/2015", 2) | Out-Null

# All eight "SHORT" title placeholders (one per section title page) become "UH HILO".
$d.Content.Find.Execute("SHORT", $true, $false, $false, $false, $false, $true, 1, $false, "UH HILO", 2) | Out-Null

# --- Move the "_GoBack" bookmark to the very first paragraph ----------
# Word keeps a single "_GoBack" bookmark marking the last edit location;
# re-adding it at the new spot automatically removes the old one.
$start = $d.Paragraphs.First.Range.Start
$r = $d.Range($start, $start)
$d.Bookmarks.Add("_GoBack", $r) | Out-Null
